# Implementando a aba da UF
# Add a "Codigo" column (C) with the numeric UF code for every state row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Header
$ws.Range("C1").Value = "Codigo"

# State -> code values, row by row (rows 2..28, matching the existing
# Estado/Sigla rows already present in column A/B)
$codigos = @(12, 27, 16, 13, 29, 23, 53, 32, 52, 21, 51, 50, 31, 15, 25, 41, 26, 22, 33, 24, 43, 11, 14, 42, 35, 28, 17)

for ($i = 0; $i -lt $codigos.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $codigos[$i]
}

# Page margins (narrow / metric defaults: 1.3cm sides, 2cm top/bottom, 0.8cm header/footer)
$ps = $ws.PageSetup
$ps.LeftMargin = 36.850393728
$ps.RightMargin = 36.850393728
$ps.TopMargin = 56.692913399999995
$ps.BottomMargin = 56.692913399999995
$ps.HeaderMargin = 22.67716464
$ps.FooterMargin = 22.67716464

# Restore the view so the previously-used cell is selected again
$ws.Range("C22").Select() | Out-Null
